$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 2 (2014/12 IFRS연결)
# ---------------------------------------------------------------------------
$ws.Range("D2").Value  = 13459
$ws.Range("E2").Value  = 439
$ws.Range("F2").Value  = 439
$ws.Range("G2").Value  = 134
$ws.Range("H2").Value  = 65
$ws.Range("I2").Value  = 65
$ws.Range("J2").Value  = 0
$ws.Range("K2").Value  = 15613
$ws.Range("L2").Value  = 11138
$ws.Range("M2").Value  = 4475
$ws.Range("N2").Value  = 4473
$ws.Range("O2").Value  = 2
$ws.Range("P2").Value  = 240
$ws.Range("Q2").Value  = -123
$ws.Range("R2").Value  = -296
$ws.Range("S2").Value  = 228
$ws.Range("T2").Value  = 171
$ws.Range("U2").Value  = -295
$ws.Range("V2").Value  = 5269
$ws.Range("W2").Value  = 3.26
$ws.Range("X2").Value  = 0.49
$ws.Range("Y2").Value  = 1.46
$ws.Range("Z2").Value  = 0.42
$ws.Range("AA2").Value = 248.9
$ws.Range("AB2").Value = 1336.46
$ws.Range("AC2").Value = 678
$ws.Range("AD2").Value = 36.82
$ws.Range("AE2").Value = 51033
$ws.Range("AF2").Value = 0.49
$ws.Range("AG2").Value = 550
$ws.Range("AH2").Value = 2.2
$ws.Range("AI2").Value = 74.28
$ws.Range("AJ2").Value = 8840002

# ---------------------------------------------------------------------------
# Row 3 (2015/12 IFRS연결)
# ---------------------------------------------------------------------------
$ws.Range("D3").Value  = 12458
$ws.Range("E3").Value  = 229
$ws.Range("F3").Value  = 303
$ws.Range("G3").Value  = -166
$ws.Range("H3").Value  = -176
$ws.Range("I3").Value  = -176
$ws.Range("J3").Value  = 1
$ws.Range("K3").Value  = 15417
$ws.Range("L3").Value  = 11168
$ws.Range("M3").Value  = 4249
$ws.Range("N3").Value  = 4246
$ws.Range("O3").Value  = 2
$ws.Range("P3").Value  = 240
$ws.Range("Q3").Value  = -73
$ws.Range("R3").Value  = -94
$ws.Range("S3").Value  = 95
$ws.Range("T3").Value  = 134
$ws.Range("U3").Value  = -208
$ws.Range("V3").Value  = 5258
$ws.Range("W3").Value  = 1.84
$ws.Range("X3").Value  = -1.41
$ws.Range("Y3").Value  = -4.05
$ws.Range("Z3").Value  = -1.13
$ws.Range("AA3").Value = 262.86
$ws.Range("AB3").Value = 1240.87
$ws.Range("AC3").Value = -1834
$ws.Range("AD3").Value = -10.99
$ws.Range("AE3").Value = 48442
$ws.Range("AF3").Value = 0.42
$ws.Range("AG3").Value = 450
$ws.Range("AH3").Value = 2.23
$ws.Range("AI3").Value = -22.47
$ws.Range("AJ3").Value = 8840002

# ---------------------------------------------------------------------------
# Row 4 (2016/12 IFRS연결)
# ---------------------------------------------------------------------------
$ws.Range("D4").Value  = 12642
$ws.Range("E4").Value  = 448
$ws.Range("F4").Value  = 448
$ws.Range("G4").Value  = 363
$ws.Range("H4").Value  = 226
$ws.Range("I4").Value  = 225
$ws.Range("J4").Value  = 1
$ws.Range("K4").Value  = 12106
$ws.Range("L4").Value  = 7608
$ws.Range("M4").Value  = 4498
$ws.Range("N4").Value  = 4495
$ws.Range("O4").Value  = 3
$ws.Range("P4").Value  = 240
$ws.Range("Q4").Value  = 869
$ws.Range("R4").Value  = 613
$ws.Range("S4").Value  = -802
$ws.Range("T4").Value  = 117
$ws.Range("U4").Value  = 753
$ws.Range("V4").Value  = 4802
$ws.Range("W4").Value  = 3.54
$ws.Range("X4").Value  = 1.79
$ws.Range("Y4").Value  = 5.16
$ws.Range("Z4").Value  = 1.64
$ws.Range("AA4").Value = 169.12
$ws.Range("AB4").Value = 1324.13
$ws.Range("AC4").Value = 2344
$ws.Range("AD4").Value = 10.24
$ws.Range("AE4").Value = 49752
$ws.Range("AF4").Value = 0.48
$ws.Range("AG4").Value = 650
$ws.Range("AH4").Value = 2.71
$ws.Range("AI4").Value = 26.14
$ws.Range("AJ4").Value = 8840002

# ---------------------------------------------------------------------------
# Row 5 (2017/12 IFRS연결)
# ---------------------------------------------------------------------------
$ws.Range("D5").Value  = 12745
$ws.Range("E5").Value  = 164
$ws.Range("F5").Value  = 164
$ws.Range("G5").Value  = 140
$ws.Range("H5").Value  = 108
$ws.Range("I5").Value  = 106
$ws.Range("J5").Value  = 2
$ws.Range("K5").Value  = 11529
$ws.Range("L5").Value  = 7011
$ws.Range("M5").Value  = 4518
$ws.Range("N5").Value  = 4513
$ws.Range("O5").Value  = 5
$ws.Range("P5").Value  = 240
$ws.Range("Q5").Value  = 105
$ws.Range("R5").Value  = -132
$ws.Range("S5").Value  = -246
$ws.Range("T5").Value  = 138
$ws.Range("U5").Value  = -33
$ws.Range("V5").Value  = 4561
$ws.Range("W5").Value  = 1.29
$ws.Range("X5").Value  = 0.85
$ws.Range("Y5").Value  = 2.36
$ws.Range("Z5").Value  = 0.92
$ws.Range("AA5").Value = 155.18
$ws.Range("AB5").Value = 1338.5
$ws.Range("AC5").Value = 1103
$ws.Range("AD5").Value = 22.22
$ws.Range("AE5").Value = 49947
$ws.Range("AF5").Value = 0.49
$ws.Range("AG5").Value = 550
$ws.Range("AH5").Value = 2.24
$ws.Range("AI5").Value = 47.04
$ws.Range("AJ5").Value = 8840002

# ---------------------------------------------------------------------------
# Row 6 (2018/12 IFRS연결) - note: J6 and O6 do not exist in this row
# ---------------------------------------------------------------------------
$ws.Range("D6").Value  = 12101
$ws.Range("E6").Value  = 289
$ws.Range("F6").Value  = 307
$ws.Range("G6").Value  = 149
$ws.Range("H6").Value  = 98
$ws.Range("I6").Value  = 97
$ws.Range("K6").Value  = 11354
$ws.Range("L6").Value  = 6835
$ws.Range("M6").Value  = 4520
$ws.Range("N6").Value  = 4515
$ws.Range("P6").Value  = 240
$ws.Range("Q6").Value  = 126
$ws.Range("R6").Value  = -180
$ws.Range("S6").Value  = -230
$ws.Range("T6").Value  = 97
$ws.Range("U6").Value  = 29
$ws.Range("V6").Value  = 4353
$ws.Range("W6").Value  = 2.38
$ws.Range("X6").Value  = 0.8100000000000001
$ws.Range("Y6").Value  = 2.15
$ws.Range("Z6").Value  = 0.85
$ws.Range("AA6").Value = 151.22
$ws.Range("AB6").Value = 1354.83
$ws.Range("AC6").Value = 1012
$ws.Range("AD6").Value = 18.78
$ws.Range("AE6").Value = 49966
$ws.Range("AF6").Value = 0.38
$ws.Range("AG6").Value = 550
$ws.Range("AH6").Value = 2.89
$ws.Range("AI6").Value = 51.28
$ws.Range("AJ6").Value = 8840002

# ---------------------------------------------------------------------------
# Rows 7-9 (2019/12(E), 2020/12(E), 2021/12(E)) - forecast rows removed,
# only the row header / label columns (A, B, C) remain.
# ---------------------------------------------------------------------------
$ws.Range("D7:AJ9").ClearContents()
